$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.042.38"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3
$ws.Range("D3").Value = "2.304.39"
$ws.Range("E3").Value = "  -0.77%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'300.85"

# Row 6
$ws.Range("D6").Value = "'98.48"
$ws.Range("E6").Value = "  -0.97%  "

# Row 7
$ws.Range("D7").Value = "'0.519"
$ws.Range("E7").Value = "  +2.37%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -1.16%  "

# Row 10
$ws.Range("D10").Value = "'36.24"
$ws.Range("E10").Value = "  -0.09%  "

# Row 11
$ws.Range("E11").Value = "  -0.41%  "

# Row 12
$ws.Range("D12").Value = "'17.88"
$ws.Range("E12").Value = "  +1.31%  "

# Row 13
$ws.Range("E13").Value = "  +0.80%  "

# Row 14
$ws.Range("E14").Value = "  -1.57%  "

# Row 15
$ws.Range("D15").Value = "2.663.20"
$ws.Range("E15").Value = "  -0.74%  "

# Row 16
$ws.Range("D16").Value = "2.304.50"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17
$ws.Range("D17").Value = "'0.782"
$ws.Range("E17").Value = "  -1.95%  "

# Row 18
$ws.Range("D18").Value = "43.000.47"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").Value = "'12.66"
$ws.Range("E19").Value = "  -2.32%  "

# Row 20
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("E21").Value = "  -1.65%  "

# Row 22
$ws.Range("D22").Value = "'68.32"
$ws.Range("E22").Value = "  +0.23%  "

# Row 23
$ws.Range("D23").Value = "'242.49"
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  -0.82%  "

# Row 28
$ws.Range("D28").Value = "'25.29"
$ws.Range("E28").Value = "  -1.06%  "

# Row 29
$ws.Range("D29").Value = "'167.01"
$ws.Range("E29").Value = "  -1.17%  "

# Row 30
$ws.Range("E30").Value = "  -0.02%  "

# Row 32
$ws.Range("D32").Value = "'33.15"
$ws.Range("E32").Value = "  -3.02%  "

# Row 33
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.03"
$ws.Range("E34").Value = "  -2.92%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'4.74"
$ws.Range("E35").Value = "  +0.27%  "

# Row 36
$ws.Range("D36").Value = "'17.74"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("D38").Value = "'0.0690"
$ws.Range("E38").Value = "  -1.12%  "

# Row 39
$ws.Range("E39").Value = "  -1.84%  "

# Row 40
$ws.Range("E40").Value = "  -1.47%  "

# Row 41
$ws.Range("E41").Value = "  +0.15%  "

# Row 42
$ws.Range("E42").Value = "  +0.92%  "

# Row 43
$ws.Range("D43").Value = "2.007.05"
$ws.Range("E43").Value = "  +0.22%  "

# Row 44
$ws.Range("D44").Value = "'0.0287"
$ws.Range("E44").Value = "  -1.03%  "

# Row 45
$ws.Range("D45").Value = "'2.17"
$ws.Range("E45").Value = "  -3.03%  "

# Row 46
$ws.Range("D46").Value = "'10.18"
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("D47").Value = "'17.47"
$ws.Range("E47").Value = "  -0.82%  "

# Row 48
$ws.Range("E48").Value = "  -2.33%  "

# Row 49
$ws.Range("D49").Value = "'53.61"
$ws.Range("E49").Value = "  -2.46%  "

# Row 50
$ws.Range("D50").Value = "2.528.99"
$ws.Range("E50").Value = "  -0.73%  "

# Row 51
$ws.Range("D51").Value = "'72.75"
$ws.Range("E51").Value = "  -5.17%  "
